# Auto-generated edit script applying the diff to Behemoth_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""
$ws.Range("H49").Value = 975.7143
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = ""
$ws.Range("N97").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").Value = ""
$ws.Range("H116").Value = 7014.2856
$ws.Range("I116").Value = 6620
$ws.Range("K116").Value = 6620
$ws.Range("M116").Value = -3178

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3221.9
$ws.Range("I45").Value = 3135.4443
$ws.Range("K45").Value = 3135.4443
$ws.Range("M45").Value = -2758.4443
$ws.Range("H61").Value = 23864648
$ws.Range("I61").Value = 33338910
$ws.Range("K61").Value = 33338910
$ws.Range("M61").Value = -33338698
$ws.Range("H74").Value = 5959147.5
$ws.Range("I74").Value = 8066282
$ws.Range("K74").Value = 8066282
$ws.Range("M74").Value = -8065408
$ws.Range("H77").Value = 5959147.5
$ws.Range("I77").Value = 8066282
$ws.Range("K77").Value = 40331410
$ws.Range("M77").Value = -40327042
$ws.Range("H97").Value = 1354.5
$ws.Range("I97").Value = 1198.8823
$ws.Range("K97").Value = 1198.8823
$ws.Range("M97").Value = -702.8823
$ws.Range("H110").Value = 2863.7058
$ws.Range("I110").Value = 2732.6
$ws.Range("K110").Value = 2732.6
$ws.Range("M110").Value = -687.5999999999999
$ws.Range("H122").Value = 1375.909
$ws.Range("I122").Value = 1376.5714
$ws.Range("K122").Value = 4129.7142
$ws.Range("M122").Value = -1679.7142
$ws.Range("H123").Value = 87995
$ws.Range("J123").Value = 87995
$ws.Range("L123").Value = 87995
$ws.Range("N123").Value = -97795
$ws.Range("H136").Value = 23864648
$ws.Range("I136").Value = 33338910
$ws.Range("K136").Value = 100016730
$ws.Range("M136").Value = -100014180

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 4500
$ws.Range("J44").Value = 4500
$ws.Range("L44").Value = 4500
$ws.Range("N44").Value = -5494
$ws.Range("H94").Value = 1842.6
$ws.Range("I94").Value = 1880.6666
$ws.Range("K94").Value = 1880.6666
$ws.Range("M94").Value = -1429.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 451155.8
$ws.Range("I31").Value = 3059.6875
$ws.Range("K31").Value = 3059.6875
$ws.Range("M31").Value = -2764.6875
$ws.Range("H34").Value = 451155.8
$ws.Range("I34").Value = 3059.6875
$ws.Range("K34").Value = 3059.6875
$ws.Range("M34").Value = -2857.6875
$ws.Range("H58").Value = 1289.25
$ws.Range("J58").Value = 1480.5
$ws.Range("L58").Value = 1480.5
$ws.Range("N58").Value = -1886.5
$ws.Range("H119").Value = 122000
$ws.Range("J119").Value = 122000
$ws.Range("L119").Value = 122000
$ws.Range("N119").Value = -131676
$ws.Range("H125").Value = 41440.332
$ws.Range("J125").Value = 41440.332
$ws.Range("L125").Value = 41440.332
$ws.Range("N125").Value = -46360.332
$ws.Range("H127").Value = 81247.5
$ws.Range("J127").Value = 81247.5
$ws.Range("L127").Value = 81247.5
$ws.Range("N127").Value = -91167.5
$ws.Range("H132").Value = 3600
$ws.Range("I132").Value = 3614.2856
$ws.Range("K132").Value = 10842.8568
$ws.Range("M132").Value = -8312.856800000001
$ws.Range("H136").Value = 1289.25
$ws.Range("J136").Value = 1480.5
$ws.Range("L136").Value = 4441.5
$ws.Range("N136").Value = -9541.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1630.8667
$ws.Range("I122").Value = 1079.8
$ws.Range("J122").Value = 1906.4
$ws.Range("K122").Value = 9718.199999999999
$ws.Range("L122").Value = 17157.6
$ws.Range("M122").Value = -7268.199999999999
$ws.Range("N122").Value = -22057.6
$ws.Range("H131").Value = 6691.1113
$ws.Range("J131").Value = 4367.3335
$ws.Range("L131").Value = 13102.0005
$ws.Range("N131").Value = -23182.0005
$ws.Range("H132").Value = 1428.1666
$ws.Range("J132").Value = 2494
$ws.Range("L132").Value = 22446
$ws.Range("N132").Value = -27506
$ws.Range("H134").Value = 3584.2666
$ws.Range("I134").Value = 2330.348
$ws.Range("K134").Value = 6991.044
$ws.Range("M134").Value = -1921.044
$ws.Range("H140").Value = 432553.56
$ws.Range("I140").Value = 603949.2
$ws.Range("K140").Value = 1811847.6
$ws.Range("M140").Value = -1806667.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 515.05884
$ws.Range("I97").Value = 511.23077
$ws.Range("J97").Value = 527.5
$ws.Range("K97").Value = 511.23077
$ws.Range("L97").Value = 527.5
$ws.Range("M97").Value = -15.23077000000001
$ws.Range("N97").Value = -1519.5
$ws.Range("H107").Value = 1812.2941
$ws.Range("I107").Value = 1127.1818
$ws.Range("J107").Value = 3068.3333
$ws.Range("K107").Value = 1127.1818
$ws.Range("L107").Value = 3068.3333
$ws.Range("M107").Value = 792.8181999999999
$ws.Range("N107").Value = -6908.3333
$ws.Range("H113").Value = 3667.15
$ws.Range("J113").Value = 4362.25
$ws.Range("L113").Value = 4362.25
$ws.Range("N113").Value = -8702.25
$ws.Range("H122").Value = 1887
$ws.Range("I122").Value = 1887
$ws.Range("K122").Value = 5661
$ws.Range("M122").Value = -3211
$ws.Range("H132").Value = 23258272
$ws.Range("I132").Value = 25002464
$ws.Range("K132").Value = 75007392
$ws.Range("M132").Value = -75004862

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3124.2
$ws.Range("J16").Value = 3848.25
$ws.Range("L16").Value = 3848.25
$ws.Range("N16").Value = -4188.25
$ws.Range("H22").Value = 2324.0833
$ws.Range("I22").Value = 3182.8333
$ws.Range("K22").Value = 3182.8333
$ws.Range("M22").Value = -2887.8333
$ws.Range("H27").Value = 2324.0833
$ws.Range("I27").Value = 3182.8333
$ws.Range("K27").Value = 3182.8333
$ws.Range("M27").Value = -3075.8333
$ws.Range("H46").Value = 3665.7036
$ws.Range("I46").Value = 3301.5
$ws.Range("J46").Value = 4195.4546
$ws.Range("K46").Value = 3301.5
$ws.Range("L46").Value = 4195.4546
$ws.Range("M46").Value = -3113.5
$ws.Range("N46").Value = -4571.4546
$ws.Range("H100").Value = 3884.8572
$ws.Range("I100").Value = 3899
$ws.Range("K100").Value = 3899
$ws.Range("M100").Value = -3358
$ws.Range("H132").Value = 155725.58
$ws.Range("I132").Value = 113016.336
$ws.Range("J132").Value = 232602.2
$ws.Range("K132").Value = 339049.008
$ws.Range("L132").Value = 697806.6000000001
$ws.Range("M132").Value = -336519.008
$ws.Range("N132").Value = -702866.6000000001
$ws.Range("H136").Value = 45708.18
$ws.Range("I136").Value = 5054.3
$ws.Range("K136").Value = 15162.9
$ws.Range("M136").Value = -12612.9

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 21740990
$ws.Range("I107").Value = 38463600
$ws.Range("J107").Value = 1600.1
$ws.Range("K107").Value = 115390800
$ws.Range("L107").Value = 4800.299999999999
$ws.Range("M107").Value = -115388880
$ws.Range("N107").Value = -8640.299999999999
$ws.Range("H113").Value = 724.5599999999999
$ws.Range("I113").Value = 724.5599999999999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2173.68
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -3.679999999999836
$ws.Range("H126").Value = 1664.625
$ws.Range("I126").Value = 1258.9333
$ws.Range("K126").Value = 3776.7999
$ws.Range("M126").Value = -1306.7999
$ws.Range("H132").Value = 16603.36
$ws.Range("I132").Value = 2613.1428
$ws.Range("K132").Value = 7839.428400000001
$ws.Range("M132").Value = -5309.428400000001
$ws.Range("H136").Value = 10818.379
$ws.Range("I136").Value = 1397.8235
$ws.Range("J136").Value = 24164.166
$ws.Range("K136").Value = 4193.470499999999
$ws.Range("L136").Value = 72492.49800000001
$ws.Range("M136").Value = -1643.470499999999
$ws.Range("N136").Value = -77592.49800000001
